$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.274.08"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "3.503.02"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'586.57"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "'134.44"
$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").Value = "4.098.64"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  +3.10%  "

$ws.Range("D15").Value = "3.500.30"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").Value = "'26.23"
$ws.Range("E16").Value = "  -4.25%  "

$ws.Range("D17").Value = "64.297.52"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "'9.93"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Value = "'13.75"
$ws.Range("E20").Value = "  -3.95%  "

$ws.Range("D21").Value = "'393.63"
$ws.Range("E21").Value = "  +2.76%  "

$ws.Range("D22").Value = "'0.573"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "3.641.34"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").Value = "'74.12"
$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").Value = "'1.53"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").Value = "3.522.54"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +3.61%  "

$ws.Range("D36").Value = "'23.50"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "'6.93"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'161.67"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'0.0784"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").Value = "'0.807"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").Value = "'25.27"
$ws.Range("E44").Value = "  -4.22%  "

$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.18"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.66"
$ws.Range("E47").Value = "  +1.79%  "

$ws.Range("D48").Value = "2.472.56"
$ws.Range("E48").Value = "  +2.10%  "

$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").Value = "'0.896"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("E51").Value = "  -1.02%  "
